$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-86 down to 29-87
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly data entry
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44477
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112001
$ws.Cells.Item(28, 7).Value = "Berenjena"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 600
$ws.Cells.Item(28, 11).Value = 8000
$ws.Cells.Item(28, 12).Value = 9000
$ws.Cells.Item(28, 13).Value = 8500
$ws.Cells.Item(28, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 142
$ws.Cells.Item(28, 17).Value = 60
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Apply the date style (s="2", numFmtId 165) to the new date cell, matching other rows in column D
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
